$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44300.90043391137
}
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44300.87771072917
}
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44300.85556826389
}
